$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data rows appended after the existing last row (125).
# Each row: Mercado ID, Mercado, Region, Fecha, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg o Unidades, Clasificacion

$ws.Range("A126").Value = 11
$ws.Range("B126").Value = "Vega Monumental Concepción"
$ws.Range("C126").Value = "Bíobío"
$ws.Range("D126").Value = 44911
$ws.Range("E126").Value = 8
$ws.Range("F126").Value = 100112028
$ws.Range("G126").Value = "Sandia"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Extra"
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 3500
$ws.Range("L126").Value = 3500
$ws.Range("M126").Value = 3500
$ws.Range("N126").Value = "`$/unidad"
$ws.Range("O126").Value = "Región de O'Higgins"
$ws.Range("P126").Value = 3500
$ws.Range("Q126").Value = 1
$ws.Range("R126").Value = "Hortaliza"

$ws.Range("A127").Value = 11
$ws.Range("B127").Value = "Vega Monumental Concepción"
$ws.Range("C127").Value = "Bíobío"
$ws.Range("D127").Value = 44911
$ws.Range("E127").Value = 8
$ws.Range("F127").Value = 100112028
$ws.Range("G127").Value = "Sandia"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 3000
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 3000
$ws.Range("N127").Value = "`$/unidad"
$ws.Range("O127").Value = "Región de O'Higgins"
$ws.Range("P127").Value = 3000
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"

$ws.Range("A128").Value = 11
$ws.Range("B128").Value = "Vega Monumental Concepción"
$ws.Range("C128").Value = "Bíobío"
$ws.Range("D128").Value = 44911
$ws.Range("E128").Value = 8
$ws.Range("F128").Value = 100112028
$ws.Range("G128").Value = "Sandia"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Segunda"
$ws.Range("J128").Value = 1000
$ws.Range("K128").Value = 2600
$ws.Range("L128").Value = 2600
$ws.Range("M128").Value = 2600
$ws.Range("N128").Value = "`$/unidad"
$ws.Range("O128").Value = "Región de O'Higgins"
$ws.Range("P128").Value = 2600
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"

# Match the date number format used by the existing Fecha column (style idx 2
# in the source file, numFmtId 165 "YYYY-MM-DD HH:MM:SS").
$ws.Range("D126:D128").NumberFormat = $ws.Range("D125").NumberFormat
